$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new work-log entry was recorded, pushing the blank "spacer" row and
# the three summary rows below it down by one. Insert a fresh row at the
# former spacer position (51); this shifts the spacer + summary rows to
# 52-54 (adding a new one at 55) and re-points the SUM() range used by
# the summary rows to the newly widened data range automatically.
[void]$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new time entry.
$ws.Range("A51").Value = 2014
$ws.Range("B51").Value = 5
$ws.Range("C51").Value = 6
$ws.Range("D51").Value = 0.55208333333333337
$ws.Range("E51").Value = 0.57291666666666663
$ws.Range("F51").Formula = "=(E51-D51)*24*60"
$ws.Range("G51").Formula = "=F51/60"

# Move the selection to A52 (the new spacer row), matching where the
# user's cursor ended up after entering the new row.
[void]$ws.Range("A52").Select()
